$wb = $excel.ActiveWorkbook

# --- Antonio sheet: fill in row 5 (Segunda semana) with new hours data ---
$wsAntonio = $wb.Worksheets.Item("Antonio")
$wsAntonio.Range("D5").Value = 2
$wsAntonio.Range("E5").Value = 1.5
$wsAntonio.Range("F5").Value = "Sábado y Domingo"
$wsAntonio.Range("H5").Value = "Si"

# --- Fernando sheet: fill in row 5 (Segunda semana) with new hours data ---
$wsFernando = $wb.Worksheets.Item("Fernando")
$wsFernando.Range("D5").Value = 2
$wsFernando.Range("E5").Value = 1
$wsFernando.Range("F5").Value = "Domingo"
$wsFernando.Range("H5").Value = "No, continua en Sprint"
